# Update the date line and all division problems in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-21 Sunday", "2025-12-22 Monday"),
    @("92÷9=", "29÷5="),
    @("66÷8=", "14÷4="),
    @("72÷5=", "51÷9="),
    @("12÷8=", "19÷9="),
    @("70÷2=", "97÷5="),
    @("73÷4=", "96÷5="),
    @("27÷4=", "33÷8="),
    @("90÷7=", "70÷8="),
    @("72÷6=", "84÷7="),
    @("87÷3=", "32÷7="),
    @("80÷3=", "39÷9="),
    @("37÷4=", "78÷4="),
    @("71÷6=", "75÷9="),
    @("79÷6=", "31÷8="),
    @("63÷9=", "45÷6="),
    @("18÷8=", "94÷9="),
    @("73÷5=", "52÷8="),
    @("23÷5=", "51÷8="),
    @("18÷2=", "61÷3="),
    @("10÷2=", "84÷2="),
    @("53÷9=", "97÷3="),
    @("14÷6=", "70÷5="),
    @("88÷6=", "35÷3="),
    @("96÷8=", "38÷2="),
    @("61÷5=", "52÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
